# "get rid of the test data"
#
# The sheet "client_persons" holds one header row (row 1) followed by
# five rows of example/test data (rows 2-6). The commit removes that
# test data: columns A-D lose their content entirely, while columns
# E-G are cleared of content but keep their existing formatting
# (style index carried by E2:G6). The header row stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the example rows, leaving the header (row 1) intact.
$ws.Range("A2:D6").ClearContents() | Out-Null
$ws.Range("E2:G6").ClearContents() | Out-Null

# Reflect the selection left behind in the sheet after the cleanup.
$ws.Range("A2:G6").Select() | Out-Null
